$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-6: date serial 45243 -> 45244
# (2023-11-13 -> 2023-11-14)
$ws.Range("C2:C6").Value = 45244
